$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (3) : "Förändrad" date stamp -> bump from 45184 to 45186 for every data row (2-46)
for ($r = 2; $r -le 46; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45186
    }
}

# Columns S, T, V, W, X, Y (19, 20, 22, 23, 24, 25): HYPERLINK(...) formulas that only had
# the URL argument now also get a friendly display-text argument equal to column A's id
# (e.g. "A 30721-2019"). Only rows 2-15 carry these hyperlink formulas.
$linkCols = 19, 20, 22, 23, 24, 25

for ($r = 2; $r -le 15; $r++) {
    $id = $ws.Cells.Item($r, 1).Value2
    foreach ($c in $linkCols) {
        $cell = $ws.Cells.Item($r, $c)
        $f = $cell.Formula
        if ($f -match '^=HYPERLINK\("([^"]*)"\)$') {
            $url = $matches[1]
            $cell.Formula = '=HYPERLINK("' + $url + '", "' + $id + '")'
        }
    }
}
